# =====================================================================
# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" worksheet right after "总计", carrying the
# fund-holdings detail for the quarter, and updates the "总计" (totals)
# sheet so the new quarter appears as its first data row (existing rows
# shift down).
# =====================================================================

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet right after "总计".
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $totalSheet)
$ws2.Name = "2022-Q4"

# Same layout as the other quarterly sheets (fund holdings table).

# Format the text columns as Text first so numeric-looking values (fund
# codes with leading zeros, "x.x0"-style percentages, etc.) are stored as
# literal strings rather than being coerced into numbers.
$ws2.Range("B1:G18").NumberFormat = "@"
$ws2.Range("H1").NumberFormat = "@"

# ---- Header row + data rows: text cells ----
$ws2.Cells.Item(1,2).Value = "基金代码"
$ws2.Cells.Item(1,3).Value = "基金名称"
$ws2.Cells.Item(1,4).Value = "基金规模"
$ws2.Cells.Item(1,5).Value = "股票总仓位"
$ws2.Cells.Item(1,6).Value = "仓位占比"
$ws2.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws2.Cells.Item(1,8).Value = "仓位排名"
$ws2.Cells.Item(2,2).Value = "003567"
$ws2.Cells.Item(2,3).Value = "华夏行业景气混合"
$ws2.Cells.Item(2,4).Value = "109.60"
$ws2.Cells.Item(2,5).Value = "93.65"
$ws2.Cells.Item(2,6).Value = "1.93"
$ws2.Cells.Item(2,7).Value = "2.1153"
$ws2.Cells.Item(3,2).Value = "010699"
$ws2.Cells.Item(3,3).Value = "东方红创新趋势混合"
$ws2.Cells.Item(3,4).Value = "24.58"
$ws2.Cells.Item(3,5).Value = "89.62"
$ws2.Cells.Item(3,6).Value = "3.62"
$ws2.Cells.Item(3,7).Value = "0.8898"
$ws2.Cells.Item(4,2).Value = "003396"
$ws2.Cells.Item(4,3).Value = "东方红优享红利混合A"
$ws2.Cells.Item(4,4).Value = "14.24"
$ws2.Cells.Item(4,5).Value = "81.32"
$ws2.Cells.Item(4,6).Value = "5.33"
$ws2.Cells.Item(4,7).Value = "0.7590"
$ws2.Cells.Item(5,2).Value = "000480"
$ws2.Cells.Item(5,3).Value = "东方红新动力灵活配置混合A"
$ws2.Cells.Item(5,4).Value = "13.65"
$ws2.Cells.Item(5,5).Value = "76.54"
$ws2.Cells.Item(5,6).Value = "4.88"
$ws2.Cells.Item(5,7).Value = "0.6661"
$ws2.Cells.Item(6,2).Value = "169103"
$ws2.Cells.Item(6,3).Value = "东方红睿轩三年定开混合"
$ws2.Cells.Item(6,4).Value = "11.89"
$ws2.Cells.Item(6,5).Value = "79.78"
$ws2.Cells.Item(6,6).Value = "5.14"
$ws2.Cells.Item(6,7).Value = "0.6111"
$ws2.Cells.Item(7,2).Value = "001564"
$ws2.Cells.Item(7,3).Value = "东方红京东大数据灵活配置混合A"
$ws2.Cells.Item(7,4).Value = "9.23"
$ws2.Cells.Item(7,5).Value = "75.39"
$ws2.Cells.Item(7,6).Value = "5.18"
$ws2.Cells.Item(7,7).Value = "0.4781"
$ws2.Cells.Item(8,2).Value = "910021"
$ws2.Cells.Item(8,3).Value = "东方红启华三年持有期混合A"
$ws2.Cells.Item(8,4).Value = "4.18"
$ws2.Cells.Item(8,5).Value = "91.17"
$ws2.Cells.Item(8,6).Value = "3.60"
$ws2.Cells.Item(8,7).Value = "0.1505"
$ws2.Cells.Item(9,2).Value = "501030"
$ws2.Cells.Item(9,3).Value = "汇添富中证环境治理指数（LOF）A"
$ws2.Cells.Item(9,4).Value = "2.99"
$ws2.Cells.Item(9,5).Value = "92.14"
$ws2.Cells.Item(9,6).Value = "1.95"
$ws2.Cells.Item(9,7).Value = "0.0583"
$ws2.Cells.Item(10,2).Value = "015769"
$ws2.Cells.Item(10,3).Value = "天弘低碳经济混合A"
$ws2.Cells.Item(10,4).Value = "1.12"
$ws2.Cells.Item(10,5).Value = "86.07"
$ws2.Cells.Item(10,6).Value = "3.77"
$ws2.Cells.Item(10,7).Value = "0.0422"
$ws2.Cells.Item(11,2).Value = "015770"
$ws2.Cells.Item(11,3).Value = "天弘低碳经济混合C"
$ws2.Cells.Item(11,4).Value = "0.99"
$ws2.Cells.Item(11,5).Value = "86.07"
$ws2.Cells.Item(11,6).Value = "3.77"
$ws2.Cells.Item(11,7).Value = "0.0373"
$ws2.Cells.Item(12,2).Value = "011313"
$ws2.Cells.Item(12,3).Value = "东方红启华三年持有期混合B"
$ws2.Cells.Item(12,4).Value = "0.89"
$ws2.Cells.Item(12,5).Value = "91.17"
$ws2.Cells.Item(12,6).Value = "3.60"
$ws2.Cells.Item(12,7).Value = "0.0320"
$ws2.Cells.Item(13,2).Value = "164908"
$ws2.Cells.Item(13,3).Value = "交银施罗德中证环境治理指数（LOF）"
$ws2.Cells.Item(13,4).Value = "1.55"
$ws2.Cells.Item(13,5).Value = "93.92"
$ws2.Cells.Item(13,6).Value = "2.00"
$ws2.Cells.Item(13,7).Value = "0.0310"
$ws2.Cells.Item(14,2).Value = "501031"
$ws2.Cells.Item(14,3).Value = "汇添富中证环境治理指数（LOF）C"
$ws2.Cells.Item(14,4).Value = "1.38"
$ws2.Cells.Item(14,5).Value = "92.14"
$ws2.Cells.Item(14,6).Value = "1.95"
$ws2.Cells.Item(14,7).Value = "0.0269"
$ws2.Cells.Item(15,2).Value = "017493"
$ws2.Cells.Item(15,3).Value = "东方红新动力灵活配置混合C"
$ws2.Cells.Item(15,4).Value = "0.35"
$ws2.Cells.Item(15,5).Value = "76.54"
$ws2.Cells.Item(15,6).Value = "4.88"
$ws2.Cells.Item(15,7).Value = "0.0171"
$ws2.Cells.Item(16,2).Value = "013413"
$ws2.Cells.Item(16,3).Value = "交银施罗德中证环境治理指数（LOF）C"
$ws2.Cells.Item(16,4).Value = "0.11"
$ws2.Cells.Item(16,5).Value = "93.92"
$ws2.Cells.Item(16,6).Value = "2.00"
$ws2.Cells.Item(16,7).Value = "0.0022"
$ws2.Cells.Item(17,2).Value = "017535"
$ws2.Cells.Item(17,3).Value = "东方红京东大数据灵活配置混合C"
$ws2.Cells.Item(17,4).Value = "0.00"
$ws2.Cells.Item(17,5).Value = "75.39"
$ws2.Cells.Item(17,6).Value = "5.18"
$ws2.Cells.Item(18,2).Value = "017536"
$ws2.Cells.Item(18,3).Value = "东方红优享红利混合C"
$ws2.Cells.Item(18,4).Value = "0.00"
$ws2.Cells.Item(18,5).Value = "81.32"
$ws2.Cells.Item(18,6).Value = "5.33"

# Drop the temporary Text number format again (the values already latched
# as strings, so this does not turn them back into numbers) so the cells
# fall back to the workbook default style, matching the other sheets.
$ws2.Range("B1:G18").ClearFormats()
$ws2.Range("H1").ClearFormats()

# ---- Numeric cells: column A (row index), column H (rank, rows 2-18), and the true-zero G17/G18 ----
# (written after the text format was cleared so these land as real numbers, not text)
$ws2.Cells.Item(2,1).Value = 0
$ws2.Cells.Item(2,8).Value = 10
$ws2.Cells.Item(3,1).Value = 1
$ws2.Cells.Item(3,8).Value = 5
$ws2.Cells.Item(4,1).Value = 2
$ws2.Cells.Item(4,8).Value = 4
$ws2.Cells.Item(5,1).Value = 3
$ws2.Cells.Item(5,8).Value = 4
$ws2.Cells.Item(6,1).Value = 4
$ws2.Cells.Item(6,8).Value = 4
$ws2.Cells.Item(7,1).Value = 5
$ws2.Cells.Item(7,8).Value = 4
$ws2.Cells.Item(8,1).Value = 6
$ws2.Cells.Item(8,8).Value = 3
$ws2.Cells.Item(9,1).Value = 7
$ws2.Cells.Item(9,8).Value = 4
$ws2.Cells.Item(10,1).Value = 8
$ws2.Cells.Item(10,8).Value = 6
$ws2.Cells.Item(11,1).Value = 9
$ws2.Cells.Item(11,8).Value = 6
$ws2.Cells.Item(12,1).Value = 10
$ws2.Cells.Item(12,8).Value = 3
$ws2.Cells.Item(13,1).Value = 11
$ws2.Cells.Item(13,8).Value = 4
$ws2.Cells.Item(14,1).Value = 12
$ws2.Cells.Item(14,8).Value = 4
$ws2.Cells.Item(15,1).Value = 13
$ws2.Cells.Item(15,8).Value = 4
$ws2.Cells.Item(16,1).Value = 14
$ws2.Cells.Item(16,8).Value = 4
$ws2.Cells.Item(17,1).Value = 15
$ws2.Cells.Item(17,7).Value = 0
$ws2.Cells.Item(17,8).Value = 4
$ws2.Cells.Item(18,1).Value = 16
$ws2.Cells.Item(18,7).Value = 0
$ws2.Cells.Item(18,8).Value = 4

# ---- Formatting: reuse the existing bold/centered/bordered header style ----
# (copy format only, from cells on "总计" that already carry it, so no new
# style entries are introduced for the parts that should look like headers)
$totalSheet.Range("B1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$ws2.Range("A2:A18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 2) Update the "总计" sheet: 2022-Q4 becomes the first data row,
#    the older quarters shift down one row.
# ---------------------------------------------------------------
$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q4"
$totalSheet.Cells.Item(2,3).Value = 17
$totalSheet.Cells.Item(2,4).Value = 5.92

$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2022-Q3"
$totalSheet.Cells.Item(3,3).Value = 6
$totalSheet.Cells.Item(3,4).Value = 2.41

$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(4,2).Value = "2021-Q4"
$totalSheet.Cells.Item(4,3).Value = 3
$totalSheet.Cells.Item(4,4).Value = 0.22

# A4 needs the same index style ("A2"/"A3" already carry it).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

